$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.881.75'
$ws.Range("E2").Value = '  -0.33%  '
$ws.Range("D3").Value = '1.630.90'
$ws.Range("E3").Value = '  -0.67%  '
$origStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = $origStyle
$ws.Range("E4").Value = '  -0.14%  '
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.39'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  -0.66%  '
$ws.Range("E6").Value = '  -0.59%  '
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("E8").Value = '  -0.85%  '
$ws.Range("E9").Value = '  -0.82%  '
$ws.Range("E10").Value = '  -0.58%  '
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0882'
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = '  -0.28%  '
$ws.Range("D12").Value = '1.863.10'
$ws.Range("E12").Value = '  -0.65%  '
$ws.Range("D13").Value = '1.630.79'
$ws.Range("E13").Value = '  -0.64%  '
$ws.Range("E14").Value = '  -1.40%  '
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.561'
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = '  -1.80%  '
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.34'
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = '  -0.36%  '
$ws.Range("D17").Value = '27.889.35'
$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '229.38'
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = '  -1.67%  '
$ws.Range("E19").Value = '  +0.98%  '
$ws.Range("D20").Value = '0.0₃0719'
$ws.Range("E20").Value = '  -0.48%  '
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.999'
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = '  -0.22%  '
$ws.Range("E22").Value = '  -1.08%  '
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.07'
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = '  -4.70%  '
$ws.Range("E24").Value = '  -0.80%  '
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.26'
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = '  +0.92%  '
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.89'
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = '  -0.33%  '
$ws.Range("E27").Value = '  -0.38%  '
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.50'
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = '  -1.26%  '
$ws.Range("E30").Value = '  -0.86%  '
$ws.Range("E31").Value = '  -0.73%  '
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.42'
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = '  +0.44%  '
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.11'
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = '  +0.41%  '
$ws.Range("D34").Value = '1.388.61'
$ws.Range("E34").Value = '  -1.51%  '
$ws.Range("E35").Value = '  -0.04%  '
$ws.Range("E36").Value = '  +9.47%  '
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.34'
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = '  -0.74%  '
$ws.Range("E38").Value = '  +1.09%  '
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.558'
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = '  -1.29%  '
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.850'
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = '  -3.43%  '
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = '  -0.17%  '
$ws.Range("B42").Value = 'WEMIXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.01'
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = '  -1.19%  '
$ws.Range("E43").Value = '  -1.90%  '
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '65.62'
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = '  -2.46%  '
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.42'
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = '  -1.93%  '
$ws.Range("D46").Value = '1.771.82'
$ws.Range("E46").Value = '  -0.72%  '
$ws.Range("E47").Value = '  -2.85%  '
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '88.33'
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = '  +0.23%  '
$ws.Range("E49").Value = '  +1.30%  '
$ws.Range("E50").Value = '  -0.55%  '
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.63'
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = '  +0.17%  '
